# "using gains for all"
# Insert two new data columns (M_TotalTax, M_CorpTax) right after the
# M_POP column (column E), shifting the existing GFA/IMF/OECD columns
# two places to the right, and populate the new columns with values.
# Also corrects the UMICs M_POP figure (E6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns at F:G - everything from the old F column
# onward (GFA - Sales, GFA - Sales + Emp, IMF - Sales, ...) slides
# right to H:O.
$ws.Range("F:G").Insert()

# New column headers.
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# New column values per row (group, in row order: High Income, LICs, LMICs, Tax haven, UMICs).
$ws.Range("F2").Value = 14106286460237.92
$ws.Range("G2").Value = 1155021202746.413

$ws.Range("F3").Value = 3207987015.574299
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 734615892234.8064
$ws.Range("G4").Value = 88889835996.30263

$ws.Range("F5").Value = 558865056646.082
$ws.Range("G5").Value = 72600947639.16805

$ws.Range("F6").Value = 4579473077980.816
$ws.Range("G6").Value = 674619880691.7614

# Corrected M_POP value for UMICs.
$ws.Range("E6").Value = 2427884184.75
